# Edit script: reformat document (remove stray proofErr markers / merge runs),
# add paragraph spacing to a list item, move a page-break marker, add a new
# "Task 2" discussion plus a "hardships & lessons learned" section, and widen
# the page margins (smaller margins = more usable space).

$d = $word.ActiveDocument

# --- Step 1: "In order to run the code..." paragraph -----------------------
# Remove the gramStart/gramEnd proofing-error markers and merge the two runs
# (which only differed because of the proofErr split) into a single run.
$p = $d.Paragraphs.Item(6)
$rng = $p.Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>In order to run the code, please follow the following instructions (also outlined in the instructions text file):</w:t></w:r></w:p>'
[void]$rng.InsertXML($xml1)

# --- Step 2: "As the assignment..." paragraph ------------------------------
# Remove the gramStart/gramEnd proofing-error markers around "Both of them"
# and merge the now-adjacent runs into one; the later ", among which" and
# ": " runs are left untouched.
$p = $d.Paragraphs.Item(10)
$rng = $p.Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>As the assignment was split into 2 tasks, I divided the work into 2 different packages. Both of them have a Main class that sets everything up</w:t></w:r><w:r><w:t>, among which</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p>'
[void]$rng.InsertXML($xml2)

# --- Step 3: "they initialize the table mapper job..." list item -----------
# Add 6pt (120 twips) of space before this paragraph.
$p = $d.Paragraphs.Item(13)
$p.Format.SpaceBefore = 6

# --- Step 4: rewrite the "Afterwards, in Task 1" .. "In Task 2, ..." run ---
# of paragraphs. This: merges the spell-checked runs ("ArrayList"/"articleID")
# into their neighbours, moves w:lastRenderedPageBreak from the "Afterwards"
# paragraph onto the "The Task 1 mapper..." paragraph, merges the bookmarked
# "write the key value pair." sentence into that same paragraph, and then
# replaces the old, truncated "In Task 2, the Reducer..." paragraph with the
# full new content: a complete Task 2 write-up, a new "hardships & lessons
# learned" Heading2 section, and its two body paragraphs (the bookmark moves
# into the first of these).
$p1 = $d.Paragraphs.Item(14)
$p3 = $d.Paragraphs.Item(16)
$rng = $d.Range($p1.Range.Start, $p3.Range.End)
$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Afterwards, in Task 1</w:t></w:r><w:r><w:t xml:space="preserve"> we have</w:t></w:r><w:r><w:t xml:space="preserve"> a custom mapper and a custom reducer. The reducer is rather simply impl</w:t></w:r><w:r><w:t>emented, adding all revisions into an ArrayList from which, after it’s sorted, we create a string containing all revision IDs in a sorted fashion. In the end, we write the key, the value’s size and the actual value pair where the key is the articleID, the size is the number of different revisions we’ve found and the value is the list above-mentioned.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>The Task 1 mapper first gets the 2 timestamps from the command line as the last 2 arguments. Then, it checks if the current timestamp is between the 2 provided by the user and, if so, we write the key value pair.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In Task 2, the </w:t></w:r><w:r><w:t>reducer</w:t></w:r><w:r><w:t xml:space="preserve"> has been developed such that it performs as optimally as possible. As it loops through t</w:t></w:r><w:r><w:t>he revisions we have, eventually it will select the one that is the closest to the timestamp provided as an argument. In the end, it is formatted from milliseconds to the ISO8601 format requested by the spec sheet and the key value pair is written. I’ve decided to use a helper class (i.e. UtilityPairRevisionTimestamp) that is in fact a pair formed of a timestamp and a revisionID. This made it easier to set the output class and manipulate the data wherever needed.</w:t></w:r></w:p><w:p><w:r><w:t>Task 2 mapper gets the timestamp from the command line arguments and checks if the timestamp residing in value is before the one passed by the user. If so, we create a new utility pair object and write the key, composite value pair.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>hardships &amp; lessons learned</w:t></w:r></w:p><w:p><w:r><w:t>The biggest hardship in my opinion was the implementation of the utility class</w:t></w:r><w:r><w:t xml:space="preserve"> for task two</w:t></w:r><w:r><w:t xml:space="preserve">. I was not sure how to proceed with getting an easy to use </w:t></w:r><w:r><w:t>data structure that can store the value, but in the end I managed to come up with the current solution. Another hardship was, without doubt, the huge latency from Hadoop – if there were even a couple of students running jobs, the output would be produced extremely slow. I had to wait until times when the load was lighter as well so that I could te</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">st my programs in all possible conditions. Another slightly hard task was to be sure that the revision IDs are sorted, but in the end this was also solved by implementing the necessary </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">As in task one, I learned that the network latency was the biggest problem in evaluating efficiency. Moreover, I learned </w:t></w:r></w:p>'
[void]$rng.InsertXML($xml4)

# --- Step 5: page margins ---------------------------------------------------
# top/right/bottom/left all become 720 twips (36pt / 0.5in).
$ps = $d.Sections.Item(1).PageSetup
$ps.TopMargin = 36
$ps.BottomMargin = 36
$ps.LeftMargin = 36
$ps.RightMargin = 36
